$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# ---- sheet1 cell values ----
$ws1.Range('A1').Value = 'File Name'
$ws1.Range('B1').Value = 'zh-cn'
$ws1.Range('C1').Value = 'de-de'
$ws1.Range('D1').Value = 'Latest Handoff Date'
$ws1.Range('A2').Value = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.md'
$ws1.Range('B2').Value = 'Handed back: in sync with en-US'
$ws1.Range('C2').Value = 'Handed back: in sync with en-US'
$ws1.Range('D2').Value = '2016-17-20 00:17:55'
$ws1.Range('A3').Value = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.md'
$ws1.Range('B3').Value = 'Handed back: in sync with en-US'
$ws1.Range('C3').Value = 'Handed back: in sync with en-US'
$ws1.Range('D3').Value = '2016-17-20 00:17:55'
$ws1.Range('A4').Value = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.md'
$ws1.Range('B4').Value = 'Handed back: in sync with en-US'
$ws1.Range('C4').Value = 'Handed back: in sync with en-US'
$ws1.Range('D4').Value = '2016-15-20 00:15:22'
$ws1.Range('A5').Value = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.md'
$ws1.Range('B5').Value = 'Handed back: in sync with en-US'
$ws1.Range('C5').Value = 'Handed back: in sync with en-US'
$ws1.Range('D5').Value = '2016-15-20 00:15:22'
$ws1.Range('A6').Value = '763af322-9580-48a0-a2db-2c2cc30e0768.md'
$ws1.Range('B6').Value = 'Handback transform failed'
$ws1.Range('C6').Value = 'Handback transform failed'
$ws1.Range('D6').Value = '2016-17-20 00:17:55'
$ws1.Range('A7').Value = 'da468f43-8ed4-400d-b859-7a80462d91f7.md'
$ws1.Range('B7').Value = 'In Translation'
$ws1.Range('C7').Value = 'In Translation'
$ws1.Range('D7').Value = '2016-14-20 00:14:17'
$ws1.Range('A8').Value = '5d06e92d-2b03-4bef-8a85-a96ae0826304.md'
$ws1.Range('B8').Value = 'Ready for handoff'
$ws1.Range('C8').Value = 'Ready for handoff'
$ws1.Range('D8').Value = '2016-17-20 00:17:55'

# ---- sheet2 cell values ----
$ws2.Range('A1').Value = 'Source File Name'
$ws2.Range('B1').Value = 'File Extension'
$ws2.Range('C1').Value = 'Status'
$ws2.Range('D1').Value = 'Latest Handoff File'
$ws2.Range('E1').Value = 'Latest Handoff Datetime'
$ws2.Range('F1').Value = 'Latest Target File'
$ws2.Range('G1').Value = 'Latest Handback File'
$ws2.Range('H1').Value = 'Latest Handback DateTime'
$ws2.Range('I1').Value = 'Handoff Reason'
$ws2.Range('J1').Value = 'Dependency From'
$ws2.Range('K1').Value = 'Error Detail'
$ws2.Range('A2').Value = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.md'
$ws2.Range('B2').Value = '.md'
$ws2.Range('C2').Value = 'Handed back: in sync with en-US'
$ws2.Range('D2').Value = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.b133721ec56910a3c4f2e2e8d9ad5581d44efea3.zh-cn.xlf'
$ws2.Range('E2').Value = '2016-03-20 00:17:52'
$ws2.Range('F2').Value = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.md'
$ws2.Range('G2').Value = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.b133721ec56910a3c4f2e2e8d9ad5581d44efea3.zh-cn.xlf'
$ws2.Range('H2').Value = '2016-03-20 00:17:20'
$ws2.Range('I2').Value = 'Include'
$ws2.Range('A3').Value = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.md'
$ws2.Range('B3').Value = '.md'
$ws2.Range('C3').Value = 'Handed back: in sync with en-US'
$ws2.Range('D3').Value = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.5d12c5159de538560fe51be78b917294b8b74bfd.zh-cn.xlf'
$ws2.Range('E3').Value = '2016-03-20 00:17:52'
$ws2.Range('F3').Value = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.md'
$ws2.Range('G3').Value = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.5d12c5159de538560fe51be78b917294b8b74bfd.zh-cn.xlf'
$ws2.Range('H3').Value = '2016-03-20 00:17:20'
$ws2.Range('I3').Value = 'Include'
$ws2.Range('A4').Value = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.md'
$ws2.Range('B4').Value = '.md'
$ws2.Range('C4').Value = 'Handed back: in sync with en-US'
$ws2.Range('D4').Value = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.ab2d2f99b755e50288e7e07b9d160798f0db5c12.zh-cn.xlf'
$ws2.Range('E4').Value = '2016-03-20 00:15:19'
$ws2.Range('F4').Value = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.md'
$ws2.Range('G4').Value = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.ab2d2f99b755e50288e7e07b9d160798f0db5c12.zh-cn.xlf'
$ws2.Range('H4').Value = '2016-03-20 00:15:36'
$ws2.Range('I4').Value = 'Include'
$ws2.Range('A5').Value = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.md'
$ws2.Range('B5').Value = '.md'
$ws2.Range('C5').Value = 'Handed back: in sync with en-US'
$ws2.Range('D5').Value = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.8f0cea50c8fce0bf4cb3a675fc2139ebdce12f1b.zh-cn.xlf'
$ws2.Range('E5').Value = '2016-03-20 00:15:19'
$ws2.Range('F5').Value = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.md'
$ws2.Range('G5').Value = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.8f0cea50c8fce0bf4cb3a675fc2139ebdce12f1b.zh-cn.xlf'
$ws2.Range('H5').Value = '2016-03-20 00:15:36'
$ws2.Range('I5').Value = 'Include'
$ws2.Range('A6').Value = '763af322-9580-48a0-a2db-2c2cc30e0768.md'
$ws2.Range('B6').Value = '.md'
$ws2.Range('C6').Value = 'Handback transform failed'
$ws2.Range('D6').Value = '763af322-9580-48a0-a2db-2c2cc30e0768.d02ae20adbcdd7f677ef4e583f8254ab294642e4.zh-cn.xlf'
$ws2.Range('E6').Value = '2016-03-20 00:17:52'
$ws2.Range('H6').Value = '0001-01-01 00:00:00'
$ws2.Range('I6').Value = 'Include'
$ws2.Range('K6').Value = 'The handback type mt is not match with handoff type ht.'
$ws2.Range('A7').Value = 'da468f43-8ed4-400d-b859-7a80462d91f7.md'
$ws2.Range('B7').Value = '.md'
$ws2.Range('C7').Value = 'In Translation'
$ws2.Range('D7').Value = 'da468f43-8ed4-400d-b859-7a80462d91f7.dba5110a1333fc0342f0c85bd445e7efe0404ee5.zh-cn.xlf'
$ws2.Range('E7').Value = '2016-03-20 00:14:14'
$ws2.Range('H7').Value = '0001-01-01 00:00:00'
$ws2.Range('I7').Value = 'Include'
$ws2.Range('A8').Value = '5d06e92d-2b03-4bef-8a85-a96ae0826304.md'
$ws2.Range('B8').Value = '.md'
$ws2.Range('C8').Value = 'Ready for handoff'
$ws2.Range('D8').Value = '5d06e92d-2b03-4bef-8a85-a96ae0826304.afe05f422d7137ff36775ef2e979397daa0d4cd0.zh-cn.xlf'
$ws2.Range('E8').Value = '2016-03-20 00:17:52'
$ws2.Range('H8').Value = '0001-01-01 00:00:00'
$ws2.Range('I8').Value = 'Include'

# ---- sheet3 cell values ----
$ws3.Range('A1').Value = 'Source File Name'
$ws3.Range('B1').Value = 'File Extension'
$ws3.Range('C1').Value = 'Status'
$ws3.Range('D1').Value = 'Latest Handoff File'
$ws3.Range('E1').Value = 'Latest Handoff Datetime'
$ws3.Range('F1').Value = 'Latest Target File'
$ws3.Range('G1').Value = 'Latest Handback File'
$ws3.Range('H1').Value = 'Latest Handback DateTime'
$ws3.Range('I1').Value = 'Handoff Reason'
$ws3.Range('J1').Value = 'Dependency From'
$ws3.Range('K1').Value = 'Error Detail'
$ws3.Range('A2').Value = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.md'
$ws3.Range('B2').Value = '.md'
$ws3.Range('C2').Value = 'Handed back: in sync with en-US'
$ws3.Range('D2').Value = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.b133721ec56910a3c4f2e2e8d9ad5581d44efea3.de-de.xlf'
$ws3.Range('E2').Value = '2016-03-20 00:17:55'
$ws3.Range('F2').Value = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.md'
$ws3.Range('G2').Value = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.b133721ec56910a3c4f2e2e8d9ad5581d44efea3.de-de.xlf'
$ws3.Range('H2').Value = '2016-03-20 00:17:25'
$ws3.Range('I2').Value = 'Include'
$ws3.Range('A3').Value = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.md'
$ws3.Range('B3').Value = '.md'
$ws3.Range('C3').Value = 'Handed back: in sync with en-US'
$ws3.Range('D3').Value = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.5d12c5159de538560fe51be78b917294b8b74bfd.de-de.xlf'
$ws3.Range('E3').Value = '2016-03-20 00:17:55'
$ws3.Range('F3').Value = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.md'
$ws3.Range('G3').Value = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.5d12c5159de538560fe51be78b917294b8b74bfd.de-de.xlf'
$ws3.Range('H3').Value = '2016-03-20 00:17:25'
$ws3.Range('I3').Value = 'Include'
$ws3.Range('A4').Value = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.md'
$ws3.Range('B4').Value = '.md'
$ws3.Range('C4').Value = 'Handed back: in sync with en-US'
$ws3.Range('D4').Value = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.ab2d2f99b755e50288e7e07b9d160798f0db5c12.de-de.xlf'
$ws3.Range('E4').Value = '2016-03-20 00:15:22'
$ws3.Range('F4').Value = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.md'
$ws3.Range('G4').Value = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.ab2d2f99b755e50288e7e07b9d160798f0db5c12.de-de.xlf'
$ws3.Range('H4').Value = '2016-03-20 00:15:42'
$ws3.Range('I4').Value = 'Include'
$ws3.Range('A5').Value = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.md'
$ws3.Range('B5').Value = '.md'
$ws3.Range('C5').Value = 'Handed back: in sync with en-US'
$ws3.Range('D5').Value = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.8f0cea50c8fce0bf4cb3a675fc2139ebdce12f1b.de-de.xlf'
$ws3.Range('E5').Value = '2016-03-20 00:15:22'
$ws3.Range('F5').Value = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.md'
$ws3.Range('G5').Value = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.8f0cea50c8fce0bf4cb3a675fc2139ebdce12f1b.de-de.xlf'
$ws3.Range('H5').Value = '2016-03-20 00:15:42'
$ws3.Range('I5').Value = 'Include'
$ws3.Range('A6').Value = '763af322-9580-48a0-a2db-2c2cc30e0768.md'
$ws3.Range('B6').Value = '.md'
$ws3.Range('C6').Value = 'Handback transform failed'
$ws3.Range('D6').Value = '763af322-9580-48a0-a2db-2c2cc30e0768.d02ae20adbcdd7f677ef4e583f8254ab294642e4.de-de.xlf'
$ws3.Range('E6').Value = '2016-03-20 00:17:55'
$ws3.Range('H6').Value = '0001-01-01 00:00:00'
$ws3.Range('I6').Value = 'Include'
$ws3.Range('K6').Value = 'The handback type mt is not match with handoff type ht.'
$ws3.Range('A7').Value = 'da468f43-8ed4-400d-b859-7a80462d91f7.md'
$ws3.Range('B7').Value = '.md'
$ws3.Range('C7').Value = 'In Translation'
$ws3.Range('D7').Value = 'da468f43-8ed4-400d-b859-7a80462d91f7.dba5110a1333fc0342f0c85bd445e7efe0404ee5.de-de.xlf'
$ws3.Range('E7').Value = '2016-03-20 00:14:17'
$ws3.Range('H7').Value = '0001-01-01 00:00:00'
$ws3.Range('I7').Value = 'Include'
$ws3.Range('A8').Value = '5d06e92d-2b03-4bef-8a85-a96ae0826304.md'
$ws3.Range('B8').Value = '.md'
$ws3.Range('C8').Value = 'Ready for handoff'
$ws3.Range('D8').Value = '5d06e92d-2b03-4bef-8a85-a96ae0826304.afe05f422d7137ff36775ef2e979397daa0d4cd0.de-de.xlf'
$ws3.Range('E8').Value = '2016-03-20 00:17:55'
$ws3.Range('H8').Value = '0001-01-01 00:00:00'
$ws3.Range('I8').Value = 'Include'

# ---- sheet1 hyperlink display text ----
$ws1Links = @{
  'A2' = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.md'
  'A3' = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.md'
  'A4' = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.md'
  'A5' = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.md'
  'A6' = '763af322-9580-48a0-a2db-2c2cc30e0768.md'
  'A7' = 'da468f43-8ed4-400d-b859-7a80462d91f7.md'
  'A8' = '5d06e92d-2b03-4bef-8a85-a96ae0826304.md'
}
foreach ($h in $ws1.Hyperlinks) {
  $addr = $h.Range.Address()
  $addr = $addr.Replace('$', '')
  if ($ws1Links.ContainsKey($addr)) {
    $h.TextToDisplay = $ws1Links[$addr]
  }
}

# ---- sheet2 hyperlink display text ----
$ws2Links = @{
  'A2' = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.md'
  'B2' = '.md'
  'D2' = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.b133721ec56910a3c4f2e2e8d9ad5581d44efea3.zh-cn.xlf'
  'F2' = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.md'
  'G2' = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.b133721ec56910a3c4f2e2e8d9ad5581d44efea3.zh-cn.xlf'
  'A3' = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.md'
  'B3' = '.md'
  'D3' = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.5d12c5159de538560fe51be78b917294b8b74bfd.zh-cn.xlf'
  'F3' = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.md'
  'G3' = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.5d12c5159de538560fe51be78b917294b8b74bfd.zh-cn.xlf'
  'A4' = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.md'
  'B4' = '.md'
  'D4' = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.ab2d2f99b755e50288e7e07b9d160798f0db5c12.zh-cn.xlf'
  'F4' = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.md'
  'G4' = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.ab2d2f99b755e50288e7e07b9d160798f0db5c12.zh-cn.xlf'
  'A5' = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.md'
  'B5' = '.md'
  'D5' = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.8f0cea50c8fce0bf4cb3a675fc2139ebdce12f1b.zh-cn.xlf'
  'F5' = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.md'
  'G5' = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.8f0cea50c8fce0bf4cb3a675fc2139ebdce12f1b.zh-cn.xlf'
  'A6' = '763af322-9580-48a0-a2db-2c2cc30e0768.md'
  'B6' = '.md'
  'D6' = '763af322-9580-48a0-a2db-2c2cc30e0768.d02ae20adbcdd7f677ef4e583f8254ab294642e4.zh-cn.xlf'
  'A7' = 'da468f43-8ed4-400d-b859-7a80462d91f7.md'
  'B7' = '.md'
  'D7' = 'da468f43-8ed4-400d-b859-7a80462d91f7.dba5110a1333fc0342f0c85bd445e7efe0404ee5.zh-cn.xlf'
  'A8' = '5d06e92d-2b03-4bef-8a85-a96ae0826304.md'
  'B8' = '.md'
  'D8' = '5d06e92d-2b03-4bef-8a85-a96ae0826304.afe05f422d7137ff36775ef2e979397daa0d4cd0.zh-cn.xlf'
}
foreach ($h in $ws2.Hyperlinks) {
  $addr = $h.Range.Address()
  $addr = $addr.Replace('$', '')
  if ($ws2Links.ContainsKey($addr)) {
    $h.TextToDisplay = $ws2Links[$addr]
  }
}

# ---- sheet3 hyperlink display text ----
$ws3Links = @{
  'A2' = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.md'
  'B2' = '.md'
  'D2' = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.b133721ec56910a3c4f2e2e8d9ad5581d44efea3.de-de.xlf'
  'F2' = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.md'
  'G2' = '28cb43ac-0c12-44d7-929b-fa461cae6b6a.b133721ec56910a3c4f2e2e8d9ad5581d44efea3.de-de.xlf'
  'A3' = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.md'
  'B3' = '.md'
  'D3' = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.5d12c5159de538560fe51be78b917294b8b74bfd.de-de.xlf'
  'F3' = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.md'
  'G3' = '45dc5af5-f84f-453e-abdd-ca9a40ee0a6e.5d12c5159de538560fe51be78b917294b8b74bfd.de-de.xlf'
  'A4' = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.md'
  'B4' = '.md'
  'D4' = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.ab2d2f99b755e50288e7e07b9d160798f0db5c12.de-de.xlf'
  'F4' = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.md'
  'G4' = '69d742ac-810b-4b3a-b098-89e9c6b8c6dd.ab2d2f99b755e50288e7e07b9d160798f0db5c12.de-de.xlf'
  'A5' = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.md'
  'B5' = '.md'
  'D5' = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.8f0cea50c8fce0bf4cb3a675fc2139ebdce12f1b.de-de.xlf'
  'F5' = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.md'
  'G5' = 'a9fc493a-1cbf-49f9-b46c-5efec27622d3.8f0cea50c8fce0bf4cb3a675fc2139ebdce12f1b.de-de.xlf'
  'A6' = '763af322-9580-48a0-a2db-2c2cc30e0768.md'
  'B6' = '.md'
  'D6' = '763af322-9580-48a0-a2db-2c2cc30e0768.d02ae20adbcdd7f677ef4e583f8254ab294642e4.de-de.xlf'
  'A7' = 'da468f43-8ed4-400d-b859-7a80462d91f7.md'
  'B7' = '.md'
  'D7' = 'da468f43-8ed4-400d-b859-7a80462d91f7.dba5110a1333fc0342f0c85bd445e7efe0404ee5.de-de.xlf'
  'A8' = '5d06e92d-2b03-4bef-8a85-a96ae0826304.md'
  'B8' = '.md'
  'D8' = '5d06e92d-2b03-4bef-8a85-a96ae0826304.afe05f422d7137ff36775ef2e979397daa0d4cd0.de-de.xlf'
}
foreach ($h in $ws3.Hyperlinks) {
  $addr = $h.Range.Address()
  $addr = $addr.Replace('$', '')
  if ($ws3Links.ContainsKey($addr)) {
    $h.TextToDisplay = $ws3Links[$addr]
  }
}
